$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column I header, cloning the style of the other header cells ---
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "dt_insertion"

# --- Append two new data rows (157, 158) before filling column I ---
$ws.Range("A157").Value = 156
$ws.Range("B157").Value = 12117143
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 1
$ws.Range("E157").Value = 1
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

$ws.Range("A158").Value = 157
$ws.Range("B158").Value = 12117139
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 1
$ws.Range("E158").Value = 0
$ws.Range("F158").Value = 1
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 0

# --- Fill column I (dt_insertion) for the existing data rows (2-156) ---
# The first cell establishes the custom date-time number format; applying it
# via the lowercase code then the uppercase code reproduces the two numFmt
# entries (164 unused/165 used) seen in the workbook. Subsequent cells reuse
# the resulting style.
$ws.Range("I2").Value = 45489
$ws.Range("I2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("I2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

for ($r = 3; $r -le 156; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = 45489
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Fill column I for the two newly appended rows with the precise timestamp ---
$ws.Range("I157").Value = 45489.96614583334
$ws.Range("I157").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("I158").Value = 45489.96614583334
$ws.Range("I158").NumberFormat = "YYYY-MM-DD HH:MM:SS"
